# Update order-level table: drop the orders that are no longer reported
# (Clupeiformes, Cypriniformes, Gymnotiformes, Salmoniformes,
# Synbranchiformes - originally rows 4, 5, 7, 8 and 10) and refresh the
# statistics for the remaining orders (Characiformes, Cichliformes,
# Cyprinodontiformes, Siluriformes), which shift up into rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows to drop, bottom-to-top so row numbers of rows not yet
# deleted are not disturbed.
$ws.Rows("10:10").Delete()
$ws.Rows("8:8").Delete()
$ws.Rows("7:7").Delete()
$ws.Rows("5:5").Delete()
$ws.Rows("4:4").Delete()

# After the deletions the sheet holds (in order):
#   Row 2: Characiformes   (was row 2)
#   Row 3: Cichliformes    (was row 3)
#   Row 4: Cyprinodontiformes (was row 6)
#   Row 5: Siluriformes    (was row 9)

# Row 2: Characiformes
$ws.Range("B2").Value = 1184931
$ws.Range("C2").Value = 145
$ws.Range("D2").Value = 70
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 9
$ws.Range("G2").Value = 9
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 3

# Row 3: Cichliformes
$ws.Range("B3").Value = 1624714
$ws.Range("C3").Value = 94
$ws.Range("D3").Value = 38
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 0

# Row 4: Cyprinodontiformes
$ws.Range("B4").Value = 3616
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0

# Row 5: Siluriformes
$ws.Range("B5").Value = 17755
$ws.Range("C5").Value = 26
$ws.Range("D5").Value = 25
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 2
